$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: lower-case / snake_case renames ---

# A1: "Product Code*" -> "code*" (keep red "*" run, Consolas font family)
$c = $ws.Range("A1")
$c.Value = "code*"
$len = $c.Text.Length
$c.Characters($len, 1).Font.Color = 255
$c.Characters($len, 1).Font.Name = "Consolas"

# B1: "Product Name*" -> "product_name*" (keep red "*" run, Calibri font)
$c = $ws.Range("B1")
$c.Value = "product_name*"
$len = $c.Text.Length
$c.Characters($len, 1).Font.Color = 255
$c.Characters($len, 1).Font.Name = "Calibri"

# C1: "Barcode" -> "barcode"
$ws.Range("C1").Value = "barcode"

# D1: "Brand" -> "brand"
$ws.Range("D1").Value = "brand"

# E1: "Pack Size" -> "pack_size"
$ws.Range("E1").Value = "pack_size"

# F1: "Category*" -> "category*" (keep red "*" run, Calibri font)
$c = $ws.Range("F1")
$c.Value = "category*"
$len = $c.Text.Length
$c.Characters($len, 1).Font.Color = 255
$c.Characters($len, 1).Font.Name = "Calibri"

# G1: "Unit" -> "unit"
$ws.Range("G1").Value = "unit"

# H1: "Min Stock" -> "min_stock"
$ws.Range("H1").Value = "min_stock"

# I1: "Max Stock" -> "max_stock"
$ws.Range("I1").Value = "max_stock"

# --- Row 2 sample data ---

# A2: sample barcode number 102048 -> 100001
$ws.Range("A2").Value = 100001

# B2, E2, F2, G2 keep same values (Mo cola / 300 / BEVERAGE / ml)
$ws.Range("B2").Value = "Mo cola"
$ws.Range("E2").Value = 300
$ws.Range("F2").Value = "BEVERAGE"
$ws.Range("G2").Value = "ml"

# --- Selection state: whole column A selected ---
$ws.Range("A1:A1048576").Select()

$wb.Save()
